$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.285.81'
$ws.Range('E2').Value = '  -1.58%  '
$ws.Range('D3').Value = '1.576.26'
$ws.Range('E3').Value = '  -1.00%  '
$ws.Range('E4').Value = '  -0.30%  '
$ws.Range('D5').Value = '''207.98'
$ws.Range('E5').Value = '  -0.47%  '
$ws.Range('E6').Value = '  -2.16%  '
$ws.Range('E7').Value = '  -0.29%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').Value = '''0.248'
$ws.Range('E9').Value = '  -1.54%  '
$ws.Range('E10').Value = '  +0.02%  '
$ws.Range('E11').Value = '  -0.25%  '
$ws.Range('D12').Value = '1.800.42'
$ws.Range('D13').Value = '1.575.66'
$ws.Range('E13').Value = '  -1.04%  '
$ws.Range('E14').Value = '  -1.41%  '
$ws.Range('D15').Value = '''0.520'
$ws.Range('E15').Value = '  -1.68%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '27.269.50'
$ws.Range('E16').Value = '  -1.60%  '
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').Value = '''62.56'
$ws.Range('E17').Value = '  -1.10%  '
$ws.Range('D18').Value = '''215.19'
$ws.Range('E18').Value = '  -1.37%  '
$ws.Range('D19').Value = '''7.33'
$ws.Range('E19').Value = '  -0.26%  '
$ws.Range('E20').Value = '  -0.97%  '
$ws.Range('E21').Value = '  -0.26%  '
$ws.Range('E22').Value = '  -0.39%  '
$ws.Range('D23').Value = '''9.42'
$ws.Range('E23').Value = '  -3.31%  '
$ws.Range('E24').Value = '  +1.34%  '
$ws.Range('D25').Value = '''151.32'
$ws.Range('E25').Value = '  -1.62%  '
$ws.Range('E26').Value = '  -5.15%  '
$ws.Range('E27').Value = '  -0.95%  '
$ws.Range('E28').Value = '  -1.32%  '
$ws.Range('E30').Value = '  -1.96%  '
$ws.Range('E31').Value = '  -2.25%  '
$ws.Range('E32').Value = '  -1.07%  '
$ws.Range('D33').Value = '1.406.73'
$ws.Range('E34').Value = '  -1.59%  '
$ws.Range('E35').Value = '  +1.53%  '
$ws.Range('E36').Value = '  -2.33%  '
$ws.Range('E37').Value = '  -2.98%  '
$ws.Range('E38').Value = '  -2.00%  '
$ws.Range('E40').Value = '  -2.79%  '
$ws.Range('E41').Value = '  -0.25%  '
$ws.Range('E42').Value = '  +1.88%  '
$ws.Range('D43').Value = '''1.82'
$ws.Range('E43').Value = '  +3.49%  '
$ws.Range('D44').Value = '''5.36'
$ws.Range('E44').Value = '  +1.78%  '
$ws.Range('D45').Value = '''2.18'
$ws.Range('E45').Value = '  +0.29%  '
$ws.Range('D46').Value = '''63.77'
$ws.Range('E46').Value = '  -1.08%  '
$ws.Range('D47').Value = '1.712.70'
$ws.Range('E47').Value = '  -1.04%  '
$ws.Range('D48').Value = '''86.18'
$ws.Range('E48').Value = '  +0.11%  '
$ws.Range('E49').Value = '  -1.73%  '
$ws.Range('D50').Value = '''0.0955'
$ws.Range('E50').Value = '  -1.16%  '
$ws.Range('E51').Value = '  -0.44%  '
